# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Updated timestamp on the report header (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 07:07"

# --- Reorder "Fiyi" / "Dominica" and "Islas Malvinas" / "Groenlandia" ---
# (the alphabetical-ish listing swapped these pairs; the case/numbers stay
# attached to the same row, only the country label moves)
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Updated case numbers ---

# Row 15: Pakistan
$ws.Range("B15").Value = 209337
$ws.Range("C15").Value = 2825
$ws.Range("D15").Value = 98503
$ws.Range("E15").Value = 106530
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 137
$ws.Range("H15").Value = 4304

# Row 17: Alemania
$ws.Range("D17").Value = 179100
$ws.Range("E17").Value = 7251

# Row 74: Uzbekistan
$ws.Range("B74").Value = 8298
$ws.Range("C74").Value = 76
$ws.Range("E74").Value = 2779

# Row 88: Kirguistan
$ws.Range("B88").Value = 5296
$ws.Range("C88").Value = 279
$ws.Range("D88").Value = 2370
$ws.Range("E88").Value = 2869
$ws.Range("G88").Value = 7
$ws.Range("H88").Value = 57

# Row 99: Tailandia
$ws.Range("B99").Value = 3171
$ws.Range("C99").Value = 2
$ws.Range("D99").Value = 3056
$ws.Range("E99").Value = 57
